$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell reference -> new text value. Values are applied as literal
# text (not re-parsed as numbers/dates) to match the source data feed,
# which stores prices/links/percentages as plain strings.
$updates = [ordered]@{
    'D2' = '26.612.80'
    'E2' = '  +1.03%  '
    'D3' = '1.825.29'
    'E3' = '  +1.76%  '
    'E4' = '  +0.35%  '
    'B5' = 'BNB'
    'C5' = 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'
    'D5' = '308.89'
    'E5' = '  +0.74%  '
    'B6' = 'USDC'
    'C6' = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
    'D6' = '1.007'
    'E6' = '  +0.30%  '
    'D7' = '0.4675'
    'E7' = '  +3.67%  '
    'D8' = '0.3603'
    'E8' = '  +0.19%  '
    'D9' = '0.07125'
    'E9' = '  +0.75%  '
    'D10' = '0.9016'
    'E10' = '  +2.14%  '
    'D11' = '0.07736'
    'E11' = '  -0.18%  '
    'D12' = '19.40'
    'E12' = '  -0.29%  '
    'D13' = '1.830.54'
    'E13' = '  +0.97%  '
    'D14' = '5.269'
    'E14' = '  -0.10%  '
    'D15' = '6.358'
    'E15' = '  +0.72%  '
    'D16' = '87.40'
    'E16' = '  +3.06%  '
    'E17' = '  +0.32%  '
    'D18' = '0.000008538'
    'E18' = '  +0.40%  '
    'D19' = '1.007'
    'E19' = '  +0.31%  '
    'D20' = '26.652.21'
    'E20' = '  +1.07%  '
    'D21' = '14.21'
    'E21' = '  -0.21%  '
    'D22' = '5.020'
    'E22' = '  +1.12%  '
    'D23' = '10.56'
    'E23' = '  +0.46%  '
    'D24' = '1.901'
    'E24' = '  -3.39%  '
    'D25' = '152.86'
    'E25' = '  +1.25%  '
    'E26' = '  +0.68%  '
    'D27' = '1.976'
    'E27' = '  -1.52%  '
    'D28' = '113.79'
    'E28' = '  +1.71%  '
    'D29' = '4.865'
    'E29' = '  -0.06%  '
    'D30' = '0.08809'
    'E30' = '  +1.69%  '
    'D31' = '3.139'
    'E31' = '  +2.29%  '
    'D32' = '2.814'
    'E32' = '  +3.67%  '
    'D33' = '1.161'
    'E33' = '  +5.17%  '
    'D34' = '0.7357'
    'E34' = '  +1.81%  '
    'D35' = '4.436'
    'E35' = '  -0.05%  '
    'D36' = '1.079'
    'E36' = '  +1.30%  '
    'D37' = '0.01930'
    'E37' = '  +0.06%  '
    'D38' = '0.05159'
    'E38' = '  +1.46%  '
    'D39' = '2.901'
    'E39' = '  +1.32%  '
    'D40' = '6.880'
    'E40' = '  +0.57%  '
    'D41' = '0.5043'
    'E41' = '  +0.15%  '
    'D42' = '0.1495'
    'E42' = '  -1.10%  '
    'D43' = '8.041'
    'E43' = '  +0.68%  '
    'D44' = '1.007'
    'E44' = '  +0.34%  '
    'D45' = '0.4658'
    'E45' = '  +0.97%  '
    'D46' = '10.03'
    'E46' = '  +2.13%  '
    'D47' = '97.79'
    'E47' = '  -3.21%  '
    'D48' = '1.571'
    'E48' = '  -0.20%  '
    'D49' = '0.06045'
    'E49' = '  +1.49%  '
    'D50' = '64.03'
    'E50' = '  +0.07%  '
    'D51' = '35.73'
    'E51' = '  -0.52%  '
}

foreach ($ref in $updates.Keys) {
    $cell = $ws.Range($ref)
    # Force text storage so numeric-looking strings (e.g. "308.89")
    # are not auto-coerced into numbers, then drop the temporary
    # number-format override so the cell keeps its original (default) style.
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$ref]
    $cell.ClearFormats()
}
